$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 0.5586043333333334
$ws.Range("H2").Value = 1.675813
$ws.Range("I2").Value = 0.01643366487114074
$ws.Range("J2").Value = 0.01643366487114074
$ws.Range("M2").Value = 174.1282373333333
$ws.Range("N2").Value = 522.384712
$ws.Range("O2").Value = 0.985625830323027
$ws.Range("P2").Value = 0.985625830323027
$ws.Range("Q2").Value = 97.26878793009513
$ws.Range("R2").Value = 875.4190913708561
$ws.Range("S2").Value = 0.01619744458386845
$ws.Range("T2").Value = 0.01619744458386845
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 0.5586043333333334
$ws.Range("H3").Value = 1.675813
$ws.Range("I3").Value = 0.01643366487114074
$ws.Range("J3").Value = 0.01643366487114074
$ws.Range("O3").Value = 0.003686901313133159
$ws.Range("P3").Value = 0.003686901313133159
$ws.Range("Q3").Value = 0.3638504703441111
$ws.Range("R3").Value = 3.274654233097
$ws.Range("S3").Value = 0.0000605893005929990515162
$ws.Range("T3").Value = 0.0000605893005929990379637
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 0.5586043333333334
$ws.Range("H4").Value = 1.675813
$ws.Range("I4").Value = 0.01643366487114074
$ws.Range("J4").Value = 0.01643366487114074
$ws.Range("M4").Value = 1.888095
$ws.Range("N4").Value = 5.664285
$ws.Range("O4").Value = 0.01068726836383999
$ws.Range("P4").Value = 0.01068726836383999
$ws.Range("Q4").Value = 1.054698048745
$ws.Range("R4").Value = 9.492282438705
$ws.Range("S4").Value = 0.0001756309866792911
$ws.Range("T4").Value = 0.0001756309866792911
$ws.Range("H5").Value = 63.825936
$ws.Range("I5").Value = 0.625901602571932
$ws.Range("J5").Value = 0.625901602571932
$ws.Range("M5").Value = 174.1282373333333
$ws.Range("N5").Value = 522.384712
$ws.Range("O5").Value = 0.985625830323027
$ws.Range("P5").Value = 0.985625830323027
$ws.Range("Q5").Value = 3704.632577276715
$ws.Range("R5").Value = 33341.69319549044
$ws.Range("S5").Value = 0.6169047867354737
$ws.Range("T5").Value = 0.6169047867354737
$ws.Range("H6").Value = 63.825936
$ws.Range("I6").Value = 0.625901602571932
$ws.Range("J6").Value = 0.625901602571932
$ws.Range("O6").Value = 0.003686901313133159
$ws.Range("P6").Value = 0.003686901313133159
$ws.Range("S6").Value = 0.002307637440414604
$ws.Range("T6").Value = 0.002307637440414604
$ws.Range("H7").Value = 63.825936
$ws.Range("I7").Value = 0.625901602571932
$ws.Range("J7").Value = 0.625901602571932
$ws.Range("M7").Value = 1.888095
$ws.Range("N7").Value = 5.664285
$ws.Range("O7").Value = 0.01068726836383999
$ws.Range("P7").Value = 0.01068726836383999
$ws.Range("Q7").Value = 40.16981021064
$ws.Range("R7").Value = 361.52829189576
$ws.Range("S7").Value = 0.006689178396043762
$ws.Range("T7").Value = 0.006689178396043762
$ws.Range("G8").Value = 12.157548
$ws.Range("H8").Value = 36.472644
$ws.Range("I8").Value = 0.3576647325569273
$ws.Range("J8").Value = 0.3576647325569273
$ws.Range("M8").Value = 174.1282373333333
$ws.Range("N8").Value = 522.384712
$ws.Range("O8").Value = 0.985625830323027
$ws.Range("P8").Value = 0.985625830323027
$ws.Range("Q8").Value = 2116.972403535392
$ws.Range("R8").Value = 19052.75163181853
$ws.Range("S8").Value = 0.3525235990036849
$ws.Range("T8").Value = 0.3525235990036849
$ws.Range("G9").Value = 12.157548
$ws.Range("H9").Value = 36.472644
$ws.Range("I9").Value = 0.3576647325569273
$ws.Range("J9").Value = 0.3576647325569273
$ws.Range("O9").Value = 0.003686901313133159
$ws.Range("P9").Value = 0.003686901313133159
$ws.Range("Q9").Value = 7.918895887604
$ws.Range("R9").Value = 71.27006298843601
$ws.Range("S9").Value = 0.001318674572125555
$ws.Range("T9").Value = 0.001318674572125555
$ws.Range("G10").Value = 12.157548
$ws.Range("H10").Value = 36.472644
$ws.Range("I10").Value = 0.3576647325569273
$ws.Range("J10").Value = 0.3576647325569273
$ws.Range("M10").Value = 1.888095
$ws.Range("N10").Value = 5.664285
$ws.Range("O10").Value = 0.01068726836383999
$ws.Range("P10").Value = 0.01068726836383999
$ws.Range("Q10").Value = 22.95460559106
$ws.Range("R10").Value = 206.59145031954
$ws.Range("S10").Value = 0.003822458981116942
$ws.Range("T10").Value = 0.003822458981116942
